# Fixed Import: the original import shifted the table one column to the
# right (leaving a stray 0/1 game-index column in A) and only the first
# game's scores were recorded as strings instead of numbers. This also
# removes the extra (incomplete/duplicate) second game row that had been
# imported.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray leading index column (A), shifting Date/Home/HomeScore/
# Visitor/VisitorScore (previously B:F) left into A:E.
$ws.Range("A1").EntireColumn.Delete()

# Remove the second game row (now row 3: 11/13/1869, Rutgers 0, Princeton 8).
$ws.Range("A3").EntireRow.Delete()

# The score columns should hold numeric values, not text.
$ws.Range("C2").Value = 4
$ws.Range("E2").Value = 6
